$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("Data")

# --- About sheet: B6 hyperlink cell keeps same text, just refresh it (shared-string reorder is automatic) ---
$ws1.Range("B6").Value = "https://www.bls.gov/cpi/tables/supplemental-files/historical-cpi-u-201912.pdf"

# --- Data sheet: normalize year-label spacing (double space -> single space before trailing period) ---
$ws2.Range("A6").Value = "1968............................................................................. ."
$ws2.Range("A7").Value = "1969............................................................................. ."
$ws2.Range("A8").Value = "1970............................................................................. ."
$ws2.Range("A9").Value = "1971............................................................................. ."
$ws2.Range("A10").Value = "1972............................................................................. ."
$ws2.Range("A11").Value = "1973............................................................................. ."
$ws2.Range("A12").Value = "1974............................................................................. ."
$ws2.Range("A13").Value = "1975............................................................................. ."
$ws2.Range("A14").Value = "1976............................................................................. ."
$ws2.Range("A15").Value = "1977............................................................................. ."
$ws2.Range("A16").Value = "1978............................................................................. ."
$ws2.Range("A17").Value = "1979............................................................................. ."
$ws2.Range("A18").Value = "1980............................................................................. ."
$ws2.Range("A19").Value = "1981............................................................................. ."
$ws2.Range("A20").Value = "1982............................................................................. ."
$ws2.Range("A21").Value = "1983............................................................................. ."
$ws2.Range("A22").Value = "1984............................................................................. ."
$ws2.Range("A23").Value = "1985............................................................................. ."
$ws2.Range("A24").Value = "1986............................................................................. ."
$ws2.Range("A25").Value = "1987............................................................................. ."
$ws2.Range("A26").Value = "1988............................................................................. ."
$ws2.Range("A27").Value = "1989............................................................................. ."
$ws2.Range("A28").Value = "1990............................................................................. ."
$ws2.Range("A29").Value = "1991............................................................................. ."
$ws2.Range("A30").Value = "1992............................................................................. ."
$ws2.Range("A31").Value = "1993............................................................................. ."
$ws2.Range("A32").Value = "1994............................................................................. ."
$ws2.Range("A33").Value = "1995............................................................................. ."
$ws2.Range("A34").Value = "1996............................................................................. ."
$ws2.Range("A35").Value = "1997............................................................................. ."
$ws2.Range("A36").Value = "1998............................................................................. ."
$ws2.Range("A37").Value = "1999............................................................................. ."
$ws2.Range("A38").Value = "2000............................................................................. ."
$ws2.Range("A39").Value = "2001............................................................................. ."
$ws2.Range("A40").Value = "2002............................................................................. ."
$ws2.Range("A41").Value = "2003............................................................................. ."
$ws2.Range("A42").Value = "2004............................................................................. ."
$ws2.Range("A43").Value = "2005............................................................................. ."
$ws2.Range("A44").Value = "2006............................................................................. ."
$ws2.Range("A45").Value = "2007............................................................................. ."
$ws2.Range("A46").Value = "2008............................................................................. ."
$ws2.Range("A47").Value = "2009............................................................................. ."
$ws2.Range("A48").Value = "2010............................................................................. ."
$ws2.Range("A49").Value = "2011............................................................................. ."
$ws2.Range("A50").Value = "2012............................................................................. ."
$ws2.Range("A51").Value = "2013............................................................................. ."
$ws2.Range("A52").Value = "2014............................................................................. ."
$ws2.Range("A53").Value = "2015............................................................................. ."
$ws2.Range("A54").Value = "2016............................................................................. ."
$ws2.Range("A55").Value = "2017............................................................................. ."
$ws2.Range("A56").Value = "2018............................................................................. ."
$ws2.Range("A57").Value = "2019............................................................................. ."

# --- Data sheet: append new 2020 row ---
$ws2.Range("A58").Value = "2020............................................................................. ."
$ws2.Range("B58").Value = 257.55700000000002
$ws2.Range("C58").Value = 260.065
$ws2.Range("D58").Value = 258.81099999999998
$ws2.Range("E58").Value = 1.4
$ws2.Range("F58").Value = 1.2
$ws2.Range("G58").Formula = '=$D$50/D58'
$ws2.Range("G58").NumberFormat = "0.000"

# --- View state: make Data the active/selected sheet, with B58 selected; About keeps B6 selected ---
$ws1.Range("B6").Select()
$ws2.Range("B58").Select()
$ws2.Activate()

# --- Window geometry (best effort; mirrors workbookView target dims) ---
$win = $excel.ActiveWindow
try {
  $win.Left = 690
  $win.Top = 1620
  $win.Width = 12720
  $win.Height = 14955
} catch {}
